# Auto-generated PowerShell Excel COM script
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D to make room for the two newest quarters
$ws.Range("D1:E1").EntireColumn.Insert()

# Populate the new column D (most recent quarter) and column E (prior quarter)
$ws.Range("D7:E7").Style = $ws.Range("F7").Style
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8:E8").Style = $ws.Range("F8").Style
$ws.Range("D8").Value = 5357000
$ws.Range("E8").Value = 5049000
$ws.Range("D9:E9").Style = $ws.Range("F9").Style
$ws.Range("D9").Value = 4406000
$ws.Range("E9").Value = 4194000
$ws.Range("D10:E10").Style = $ws.Range("F10").Style
$ws.Range("D10").Value = 951000
$ws.Range("E10").Value = 855000
$ws.Range("D12:E12").Style = $ws.Range("F12").Style
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13:E13").Style = $ws.Range("F13").Style
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14:E14").Style = $ws.Range("F14").Style
$ws.Range("D14").Value = 91000
$ws.Range("E14").Value = 12000
$ws.Range("D15:E15").Style = $ws.Range("F15").Style
$ws.Range("D15").Value = 62000
$ws.Range("E15").Value = 52000
$ws.Range("D17:E17").Style = $ws.Range("F17").Style
$ws.Range("D17").Value = 4801000
$ws.Range("E17").Value = 4479000
$ws.Range("D18:E18").Style = $ws.Range("F18").Style
$ws.Range("D18").Value = 556000
$ws.Range("E18").Value = 570000
$ws.Range("D20:E20").Style = $ws.Range("F20").Style
$ws.Range("D20").Value = 17000
$ws.Range("E20").Value = 84000
$ws.Range("D21:E21").Style = $ws.Range("F21").Style
$ws.Range("D21").Value = 649000
$ws.Range("E21").Value = 719000
$ws.Range("D22:E22").Style = $ws.Range("F22").Style
$ws.Range("D22").Value = 94000
$ws.Range("E22").Value = 86000
$ws.Range("D23:E23").Style = $ws.Range("F23").Style
$ws.Range("D23").Value = 479000
$ws.Range("E23").Value = 568000
$ws.Range("D24:E24").Style = $ws.Range("F24").Style
$ws.Range("D24").Value = 126000
$ws.Range("E24").Value = 85000
$ws.Range("D25:E25").Style = $ws.Range("F25").Style
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26:E26").Style = $ws.Range("F26").Style
$ws.Range("D26").Value = 353000
$ws.Range("E26").Value = 483000
$ws.Range("D27:E27").Style = $ws.Range("F27").Style
$ws.Range("D27").Value = 353000
$ws.Range("E27").Value = 483000
$ws.Range("D28:E28").Style = $ws.Range("F28").Style
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29:E29").Style = $ws.Range("F29").Style
$ws.Range("D29").Value = 63000
$ws.Range("E29").Value = 0
$ws.Range("D30:E30").Style = $ws.Range("F30").Style
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31:E31").Style = $ws.Range("F31").Style
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32:E32").Style = $ws.Range("F32").Style
$ws.Range("D32").Value = -17000
$ws.Range("E32").Value = -84000
$ws.Range("D33:E33").Style = $ws.Range("F33").Style
$ws.Range("D33").Value = 416000
$ws.Range("E33").Value = 483000
$ws.Range("D34:E34").Style = $ws.Range("F34").Style
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35:E35").Style = $ws.Range("F35").Style
$ws.Range("D35").Value = 416000
$ws.Range("E35").Value = 483000
$ws.Range("D38:E38").Style = $ws.Range("F38").Style
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41:E41").Style = $ws.Range("F41").Style
$ws.Range("D41").Value = 316000
$ws.Range("E41").Value = 373000
$ws.Range("D42:E42").Style = $ws.Range("F42").Style
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43:E43").Style = $ws.Range("F43").Style
$ws.Range("D43").Value = 2133000
$ws.Range("E43").Value = 2175000
$ws.Range("D44:E44").Style = $ws.Range("F44").Style
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45:E45").Style = $ws.Range("F45").Style
$ws.Range("D45").Value = 257000
$ws.Range("E45").Value = 255000
$ws.Range("D46:E46").Style = $ws.Range("F46").Style
$ws.Range("D46").Value = 2706000
$ws.Range("E46").Value = 2803000
$ws.Range("D47:E47").Style = $ws.Range("F47").Style
$ws.Range("D47").Value = 857000
$ws.Range("E47").Value = 811000
$ws.Range("D48:E48").Style = $ws.Range("F48").Style
$ws.Range("D48").Value = 1956000
$ws.Range("E48").Value = 1967000
$ws.Range("D49:E49").Style = $ws.Range("F49").Style
$ws.Range("D49").Value = 17419000
$ws.Range("E49").Value = 17470000
$ws.Range("D50:E50").Style = $ws.Range("F50").Style
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51:E51").Style = $ws.Range("F51").Style
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52:E52").Style = $ws.Range("F52").Style
$ws.Range("D52").Value = 758000
$ws.Range("E52").Value = 782000
$ws.Range("D53:E53").Style = $ws.Range("F53").Style
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54:E54").Style = $ws.Range("F54").Style
$ws.Range("D54").Value = 23696000
$ws.Range("E54").Value = 23833000
$ws.Range("D57:E57").Style = $ws.Range("F57").Style
$ws.Range("D57").Value = 767000
$ws.Range("E57").Value = 759000
$ws.Range("D58:E58").Style = $ws.Range("F58").Style
$ws.Range("D58").Value = 833000
$ws.Range("E58").Value = 617000
$ws.Range("D59:E59").Style = $ws.Range("F59").Style
$ws.Range("D59").Value = 4837000
$ws.Range("E59").Value = 4675000
$ws.Range("D60:E60").Style = $ws.Range("F60").Style
$ws.Range("D60").Value = 6437000
$ws.Range("E60").Value = 6051000
$ws.Range("D61:E61").Style = $ws.Range("F61").Style
$ws.Range("D61").Value = 8514000
$ws.Range("E61").Value = 8710000
$ws.Range("D62:E62").Style = $ws.Range("F62").Style
$ws.Range("D62").Value = 6520000
$ws.Range("E62").Value = 6748000
$ws.Range("D63:E63").Style = $ws.Range("F63").Style
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64:E64").Style = $ws.Range("F64").Style
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65:E65").Style = $ws.Range("F65").Style
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66:E66").Style = $ws.Range("F66").Style
$ws.Range("D66").Value = 21471000
$ws.Range("E66").Value = 21509000
$ws.Range("D68:E68").Style = $ws.Range("F68").Style
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69:E69").Style = $ws.Range("F69").Style
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70:E70").Style = $ws.Range("F70").Style
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71:E71").Style = $ws.Range("F71").Style
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72:E72").Style = $ws.Range("F72").Style
$ws.Range("D72").Value = 8982000
$ws.Range("E72").Value = 8705000
$ws.Range("D73:E73").Style = $ws.Range("F73").Style
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74:E74").Style = $ws.Range("F74").Style
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75:E75").Style = $ws.Range("F75").Style
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76:E76").Style = $ws.Range("F76").Style
$ws.Range("D76").Value = 2225000
$ws.Range("E76").Value = 2324000
$ws.Range("D77:E77").Style = $ws.Range("F77").Style
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80:E80").Style = $ws.Range("F80").Style
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81:E81").Style = $ws.Range("F81").Style
$ws.Range("D81").Value = 416000
$ws.Range("E81").Value = 483000
$ws.Range("D83:E83").Style = $ws.Range("F83").Style
$ws.Range("D83").Value = 76000
$ws.Range("E83").Value = 65000
$ws.Range("D84:E84").Style = $ws.Range("F84").Style
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85:E85").Style = $ws.Range("F85").Style
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86:E86").Style = $ws.Range("F86").Style
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87:E87").Style = $ws.Range("F87").Style
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88:E88").Style = $ws.Range("F88").Style
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89:E89").Style = $ws.Range("F89").Style
$ws.Range("D89").Value = 533000
$ws.Range("E89").Value = 695000
$ws.Range("D91:E91").Style = $ws.Range("F91").Style
$ws.Range("D91").Value = -94000
$ws.Range("E91").Value = -74000
$ws.Range("D92:E92").Style = $ws.Range("F92").Style
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93:E93").Style = $ws.Range("F93").Style
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94:E94").Style = $ws.Range("F94").Style
$ws.Range("D94").Value = -133000
$ws.Range("E94").Value = -33000
$ws.Range("D96:E96").Style = $ws.Range("F96").Style
$ws.Range("D96").Value = -139000
$ws.Range("E96").Value = -142000
$ws.Range("D97:E97").Style = $ws.Range("F97").Style
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98:E98").Style = $ws.Range("F98").Style
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99:E99").Style = $ws.Range("F99").Style
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100:E100").Style = $ws.Range("F100").Style
$ws.Range("D100").Value = -453000
$ws.Range("E100").Value = -652000
$ws.Range("D101:E101").Style = $ws.Range("F101").Style
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102:E102").Style = $ws.Range("F102").Style
$ws.Range("D102").Value = -53000
$ws.Range("E102").Value = 10000

# Apply restated (revised) figures for the two quarters that were corrected in this refresh
$ws.Range("H8").Value = 5251000
$ws.Range("I8").Value = 5078000
$ws.Range("H9").Value = 4445000
$ws.Range("I9").Value = 4001000
$ws.Range("H10").Value = 806000
$ws.Range("I10").Value = 1077000
$ws.Range("H15").Value = 53000
$ws.Range("I15").Value = 54000
$ws.Range("H17").Value = 4827000
$ws.Range("I17").Value = 4288000
$ws.Range("H18").Value = 424000
$ws.Range("I18").Value = 790000
$ws.Range("H20").Value = 754000
$ws.Range("I20").Value = 94000
$ws.Range("H21").Value = 1247000
$ws.Range("I21").Value = 949000
$ws.Range("H22").Value = 144000
$ws.Range("I22").Value = 146000
$ws.Range("H23").Value = 1034000
$ws.Range("I23").Value = 738000
$ws.Range("H24").Value = 328000
$ws.Range("I24").Value = 253000
$ws.Range("H26").Value = 706000
$ws.Range("I26").Value = 485000
$ws.Range("H27").Value = 706000
$ws.Range("I27").Value = 485000
$ws.Range("H29").Value = -592000
$ws.Range("H32").Value = -754000
$ws.Range("I32").Value = -94000
$ws.Range("H33").Value = 114000
$ws.Range("I33").Value = 485000
$ws.Range("H35").Value = 114000
$ws.Range("I35").Value = 485000
$ws.Range("H81").Value = 114000
$ws.Range("I81").Value = 485000
$ws.Range("I91").Value = -51000
$ws.Range("J91").Value = -56000
$ws.Range("H94").Value = 884000
$ws.Range("I94").Value = 4000
$ws.Range("H100").Value = -1465000
$ws.Range("I100").Value = -550000
$ws.Range("H102").Value = -108000
